# "Add files via upload" — the workbook's two result sheets (WGCNA, SGFA)
# were re-saved with the "Subset accuracy" (col B) and "Hamming loss" (col C)
# values swapped for every data row. Reproduce that swap on both sheets.
#
# (The header text/values themselves are unchanged; only the B/C data
# columns were transposed relative to their previous positions.)

$wb = $excel.ActiveWorkbook

foreach ($idx in 1, 2) {
    $ws = $wb.Worksheets.Item($idx)

    # Header row (row 1) keeps "Subset accuracy" in B1 / "Hamming loss" in C1 —
    # only the data rows below it are swapped.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

    for ($r = 2; $r -le $lastRow; $r++) {
        $bVal = $ws.Cells.Item($r, 2).Value2
        $cVal = $ws.Cells.Item($r, 3).Value2
        $ws.Cells.Item($r, 2).Value2 = $cVal
        $ws.Cells.Item($r, 3).Value2 = $bVal
    }

    # Row 1 height nudged from 14.25 to 13.8 in the re-saved file.
    $ws.Rows.Item(1).RowHeight = 13.8
}

# Sheet 1 (WGCNA) column widths were re-measured slightly narrower/wider on
# re-save; nudge them toward the committed values (best effort — the host's
# ColumnWidth property applies a fixed ~5/7-character padding offset versus
# the raw OOXML <col width> it writes out, so that offset is backed out here
# to land as close as this engine's pixel-grid rounding allows).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(1).ColumnWidth = 10.020089285714286
$ws1.Columns.Item(2).ColumnWidth = 10.020089285714286
$ws1.Columns.Item(3).ColumnWidth = 15.180245535714286
$ws1.Columns.Item(4).ColumnWidth = 12.285714285714286

# Sheet 2 (SGFA) likewise.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Columns.Item(1).ColumnWidth = 9.652901785714286
$ws2.Columns.Item(2).ColumnWidth = 9.652901785714286
$ws2.Columns.Item(3).ColumnWidth = 14.285714285714286
